# Retraining the model for Horeco
# Shift all timestamps (and the derived "Lookup" text built from the date +
# quarter index) forward by 2 days, matching the new data-fetch window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -eq $null) { continue }

    $newA = $aVal + 2
    $aCell.Value = $newA

    $dt = [DateTime]::FromOADate($newA)
    $datePart = $dt.ToString("dd.MM.yyyy")

    $quarter = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = "$datePart$quarter"
}
